# Add a "Data Source" column to the NMDC-EDGE Metagenomics ReadsQC bulk
# submission template: a new column is inserted at position C (pushing the
# old C/D/E/F columns to D/E/F/G), a "Data Source" drop-down list is added
# for that new column, and the prompts/ranges for the FASTQ columns that
# moved are refreshed to reflect their new letters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Data Source" column before the old column C -------
# (Interleaved/Single FASTQ, Illumina R1 FASTQ, Illumina R2 FASTQ,
# Sequencing Platform all shift one column to the right: C->D, D->E, E->F,
# F->G.)
$ws.Range("C1").EntireColumn.Insert()

# Give the new column roughly the same width as its neighbour (column B)
# so it isn't left at the generic default width.
$ws.Range("C1").EntireColumn.ColumnWidth = $ws.Range("B1").EntireColumn.ColumnWidth()

# --- 2. Header for the new column ------------------------------------------
$ws.Range("C1").Value = "Data Source"

# --- 3. Keep the AutoFilter / _FilterDatabase defined name in sync ---------
# It used to span the 5 original header cells (A1:E1); now there are 6.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$1"
    }
}

# --- 4. Refresh the data validation prompts that reference column letters --
# (they all still apply to rows 2:99 of their respective, now-shifted,
# columns)
$newFastqPrompt = "Omit this column if input is Single fastq`n`nEnter file name if Data Source is Uploaded File or Retrieved SRA Data`n`nEnter file url if Data Source is HTTP(s) URL`n`nSeparate multiple inputs with commas"

# Illumina R1 FASTQ validation, now on column E
$r1 = $ws.Range("E2:E99").Validation
$r1.InputTitle = "Illumina R1 FASTQ"
$r1.InputMessage = $newFastqPrompt

# Single/Interleaved Illumina/PacBio FASTQ validation, now on column D
$newSinglePrompt = "Omit this column if input is Paired-end fastq`n`nEnter file name if Data Source is Uploaded File or Retrieved SRA Data`n`nEnter file url if Data Source is HTTP(s) URL`n`nSeparate multiple inputs with commas"
$single = $ws.Range("D2:D99").Validation
$single.InputTitle = "Single Illumina/PacBio FASTQ"
$single.InputMessage = $newSinglePrompt

# Illumina R2 FASTQ validation, now on column F
$r2 = $ws.Range("F2:F99").Validation
$r2.InputTitle = "Illumina R2 FASTQ"
$r2.InputMessage = $newFastqPrompt

# Sequencing Platform validation, now on column G (text unchanged)
# (no textual change needed - already shifted automatically with the
# column insert)

# --- 5. Add the new "Data Source" list validation on column C --------------
$dataSource = $ws.Range("C2:C108")
$dataSource.Validation.Add(3, 1, 1, """Uploaded File, Retrieved SRA Data, HTTP(s) URL""")
$dataSource.Validation.InputTitle = "Data Source"
$dataSource.Validation.InputMessage = "Default: Uploaded File"
$dataSource.Validation.IgnoreBlank = 1
$dataSource.Validation.InCellDropdown = 1
$dataSource.Validation.ShowInput = 1
$dataSource.Validation.ShowError = 1
